$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores values as plain text even when they look
# numeric (e.g. "583.07"). Force a Text number format on the cells we are
# about to rewrite with numeric-looking strings so Excel does not silently
# convert them to real numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"

$ws.Range("D2").Value = "67.144.96"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "2.473.69"
$ws.Range("E3").Value = "  -2.09%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "583.07"
$ws.Range("E5").Value = "  -1.53%  "
$ws.Range("D6").Value = "169.03"
$ws.Range("E6").Value = "  -1.50%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -2.06%  "
$ws.Range("D9").Value = "2.473.73"
$ws.Range("E9").Value = "  -2.01%  "
$ws.Range("E10").Value = "  -2.45%  "
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("E12").Value = "  -2.70%  "
$ws.Range("E13").Value = "  -3.21%  "
$ws.Range("D14").Value = "25.60"
$ws.Range("E14").Value = "  -3.08%  "
$ws.Range("D15").Value = "2.887.30"
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("D16").Value = "66.747.09"
$ws.Range("E16").Value = "  -1.08%  "
$ws.Range("E17").Value = "  -4.08%  "
$ws.Range("D18").Value = "2.462.17"
$ws.Range("E18").Value = "  -2.03%  "
$ws.Range("D19").Value = "11.18"
$ws.Range("E19").Value = "  -5.24%  "
$ws.Range("D20").Value = "7.61"
$ws.Range("E20").Value = "  -3.33%  "
$ws.Range("D21").Value = "354.16"
$ws.Range("E21").Value = "  -3.86%  "
$ws.Range("D22").Value = "4.04"
$ws.Range("E22").Value = "  -2.60%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "69.08"
$ws.Range("E24").Value = "  -3.54%  "
$ws.Range("E25").Value = "  -7.31%  "
$ws.Range("E26").Value = "  -6.80%  "
$ws.Range("D27").Value = "9.31"
$ws.Range("E27").Value = "  -6.56%  "
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  -1.39%  "
$ws.Range("D29").Value = "2.566.69"
$ws.Range("E29").Value = "  -2.99%  "
$ws.Range("E30").Value = "  -5.60%  "
$ws.Range("D31").Value = "517.93"
$ws.Range("E31").Value = "  -3.75%  "
$ws.Range("E32").Value = "  -7.12%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "1.78"
$ws.Range("E33").Value = "  -5.11%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "1.24"
$ws.Range("E34").Value = "  -5.68%  "
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("D36").Value = "0.119"
$ws.Range("E36").Value = "  -6.97%  "
$ws.Range("D37").Value = "157.92"
$ws.Range("E37").Value = "  -0.48%  "
$ws.Range("D38").Value = "18.67"
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("E39").Value = "  -3.33%  "
$ws.Range("E40").Value = "  -5.22%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("E42").Value = "  -6.70%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").Value = "4.80"
$ws.Range("E43").Value = "  -6.41%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "1.66"
$ws.Range("E44").Value = "  -6.30%  "
$ws.Range("E45").Value = "  -5.12%  "
$ws.Range("D46").Value = "38.71"
$ws.Range("E46").Value = "  -2.17%  "
$ws.Range("D47").Value = "141.21"
$ws.Range("E47").Value = "  -3.55%  "
$ws.Range("E48").Value = "  -6.44%  "
$ws.Range("E49").Value = "  -6.47%  "
$ws.Range("E50").Value = "  -11.14%  "
$ws.Range("E51").Value = "  -7.22%  "
